# Generate Report for Handoff
# Updates the status of the "9dceee31-3dfb-45aa-89ea-41fc7936a46f.md" file
# from "In Translation" to "Ready for handoff" across the Overview, zh-cn
# and de-de worksheets, and records the new handoff datetimes for each
# locale's handoff.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 corresponds to the 9dceee31-... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row 3 corresponds to the 9dceee31-... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-03-10 02:15:01"

# --- de-de sheet: row 3 corresponds to the 9dceee31-... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-03-10 02:15:09"
